$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price values (matches source data which
# stores these as text, e.g. using "." as a thousands separator elsewhere in the column).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = '29.199.61'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '1.829.78'
$ws.Range("E3").Value = '  -0.75%  '
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '237.75'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").Value = '0.6054'
$ws.Range("E6").Value = '  -3.88%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '0.07093'
$ws.Range("E8").Value = '  -4.61%  '
$ws.Range("D9").Value = '0.2828'
$ws.Range("E9").Value = '  -2.68%  '
$ws.Range("D11").Value = '0.07645'
$ws.Range("E11").Value = '  -1.22%  '
$ws.Range("D12").Value = '1.825.35'
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("D13").Value = '4.799'
$ws.Range("E13").Value = '  -3.73%  '
$ws.Range("D14").Value = '0.6377'
$ws.Range("E14").Value = '  -6.14%  '
$ws.Range("D15").Value = '0.000009957'
$ws.Range("E15").Value = '  -2.68%  '
$ws.Range("D16").Value = '2.066.91'
$ws.Range("E16").Value = '  -1.21%  '
$ws.Range("D17").Value = '79.74'
$ws.Range("E17").Value = '  -2.80%  '
$ws.Range("D18").Value = '5.987'
$ws.Range("E18").Value = '  -4.46%  '
$ws.Range("D19").Value = '29.176.84'
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("D20").Value = '230.46'
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("E22").Value = '  -4.39%  '
$ws.Range("D23").Value = '6.987'
$ws.Range("E23").Value = '  -5.94%  '
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = '155.57'
$ws.Range("E25").Value = '  -1.61%  '
$ws.Range("D26").Value = '8.041'
$ws.Range("E26").Value = '  -5.34%  '
$ws.Range("D27").Value = '0.1289'
$ws.Range("E27").Value = '  -4.77%  '
$ws.Range("E28").Value = '  -4.24%  '
$ws.Range("D29").Value = '0.06671'
$ws.Range("E29").Value = '  +2.20%  '
$ws.Range("D30").Value = '1.454'
$ws.Range("E30").Value = '  +0.34%  '
$ws.Range("D31").Value = '1.461'
$ws.Range("E31").Value = '  -1.73%  '
$ws.Range("D32").Value = '3.834'
$ws.Range("E32").Value = '  -5.57%  '
$ws.Range("D33").Value = '3.811'
$ws.Range("E33").Value = '  -6.41%  '
$ws.Range("E34").Value = '  -0.51%  '
$ws.Range("D35").Value = '1.718'
$ws.Range("E35").Value = '  -6.59%  '
$ws.Range("D36").Value = '0.6567'
$ws.Range("E36").Value = '  -5.69%  '
$ws.Range("E37").Value = '  -0.89%  '
$ws.Range("D38").Value = '1.238.81'
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("D39").Value = '2.755'
$ws.Range("E39").Value = '  -2.09%  '
$ws.Range("D40").Value = '0.01766'
$ws.Range("E40").Value = '  -4.70%  '
$ws.Range("D41").Value = '6.572'
$ws.Range("E41").Value = '  -3.33%  '
$ws.Range("D42").Value = '0.9283'
$ws.Range("E42").Value = '  -0.60%  '
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").Value = '1.981.78'
$ws.Range("E44").Value = '  -0.55%  '
$ws.Range("D45").Value = '100.27'
$ws.Range("D46").Value = '63.45'
$ws.Range("E46").Value = '  -3.34%  '
$ws.Range("E47").Value = '  -1.74%  '
$ws.Range("D48").Value = '1.630'
$ws.Range("E48").Value = '  -4.76%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '8.545'
$ws.Range("E49").Value = '  -5.24%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.05584'
$ws.Range("E50").Value = '  -1.49%  '
$ws.Range("D51").Value = '0.1082'
$ws.Range("E51").Value = '  -5.94%  '
